$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.078.79'
$ws.Range('D3').Value = '1.638.36'
$ws.Range('E3').Value = '  -1.77%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.96'
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  -0.81%  '
$ws.Range('E10').Value = '  -2.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07673'
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').Value = '1.638.48'
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.402'
$ws.Range('E13').Value = '  -0.56%  '
$ws.Range('D14').Value = '1.861.03'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '0.0₅8230'
$ws.Range('E16').Value = '  +3.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.88'
$ws.Range('E17').Value = '  -2.24%  '
$ws.Range('D18').Value = '26.065.29'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.682'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '188.10'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.19'
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.155'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.46'
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1211'
$ws.Range('E26').Value = '  -2.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.405'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.81'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.400'
$ws.Range('E29').Value = '  +3.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05958'
$ws.Range('E30').Value = '  -6.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.256'
$ws.Range('E31').Value = '  -1.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.431'
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.399'
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.640'
$ws.Range('E34').Value = '  +0.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9816'
$ws.Range('E35').Value = '  -1.78%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.759'
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5641'
$ws.Range('E38').Value = '  -6.19%  '
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8492'
$ws.Range('E40').Value = '  -2.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.702'
$ws.Range('E42').Value = '  -6.35%  '
$ws.Range('D43').Value = '1.031.78'
$ws.Range('E43').Value = '  -7.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.20'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').Value = '1.787.56'
$ws.Range('E45').Value = '  -1.68%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.78'
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.0000'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.025'
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05145'
$ws.Range('E49').Value = '  -1.67%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4216'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.904'
$ws.Range('E51').Value = '  -0.42%  '
